{"js": "// The paragraph originally reads \"Version 2.\" and must become \"Version 1.\".\n// Do this as two small, targeted replacements (rather than rewriting the\n// whole paragraph) so the edit stays surgical, matching how the change was\n// actually made in the source document:\n//   1) \"Version\" is re-typed as a single word (this merges the \"Versi\"/\"on\"\n//      run split left over from a spell-check correction into one run).\n//   2) \" 2.\" becomes \" 1.\" (updates the digit and folds the trailing \".\"\n//      run into the same run as the space+digit).\n\nconst body = context.document.body;\n\n// 1) Merge \"Versi\" + \"on\" into a single \"Version\" run.\nconst versionResults = body.search(\"Version\", { matchCase: true });\nversionResults.load(\"items\");\nawait context.sync();\n\nif (versionResults.items.length === 0) {\n  throw new Error('Could not find \"Version\" in the document body.');\n}\nversionResults.items[0].insertText(\"Version\", \"Replace\");\nawait context.sync();\n\n// 2) Replace \" 2.\" with \" 1.\" (bumps the version number and removes the\n// separate trailing \".\" run).\nconst suffixResults = body.search(\" 2.\", { matchCase: true });\nsuffixResults.load(\"items\");\nawait context.sync();\n\nif (suffixResults.items.length === 0) {\n  throw new Error('Could not find \" 2.\" in the document body.');\n}\nsuffixResults.items[0].insertText(\" 1.\", \"Replace\");\nawait context.sync();\n", "ps1": "# The paragraph originally reads \"Version 2.\" and must become \"Version 1.\".\n# The paragraph also contains a \"_GoBack\" bookmark sitting between the\n# \" 2\" run and the final \".\" run; a Find/Replace whose matched range spans\n# across that bookmark silently deletes it, so every edit below is kept\n# strictly on one side of the bookmark.\n\n$d = $word.ActiveDocument\n\n# 1) Re-type \"Version\" as a single word. The original run split (\"Versi\" +\n#    \"on\") is a left-over from a spell-check correction; replacing the\n#    matched range with the same text merges it into one run, same as\n#    retyping it in the UI would.\n$findVersion = $d.Content.Find\n$findVersion.Text = \"Version\"\n$findVersion.Replacement.Text = \"Version\"\n$findVersion.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 2) Delete the trailing \".\" run that sits after the bookmark (this run is\n#    entirely after the bookmark, so removing it does not disturb it).\n$findDot = $d.Content.Find\n$findDot.Text = \".\"\n$findDot.Execute()\n$findDot.Parent.Text = \"\"\n\n# 3) Rewrite the \" 2\" run (entirely before the bookmark) as \" 1.\", bumping\n#    the version number and folding the period back in as part of this run.\n$findNumber = $d.Content.Find\n$findNumber.Text = \" 2\"\n$findNumber.Execute()\n$findNumber.Parent.Text = \" 1.\"\n"}
